$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row/column, as derived from the authoritative diff.
# Columns: E, G, H, I, J, K, M, N, O, P, Q, R, S, T change; F, L stay at 1;
# A-D (text) are unchanged.

$data = @{
    2  = @{ E=3; G=1.408030333333333;  H=4.224091;  I=0.3454737251382253; J=0.3454737251382253; K=3; M=4.913477;           N=14.740431; O=0.3201889893918886; P=0.3201889893918885; Q=6.918324658135667;  R=62.264921923221;   S=0.1106168829134594;  T=0.1106168829134594 }
    3  = @{ E=3; G=1.408030333333333;  H=4.224091;  I=0.3454737251382253; J=0.3454737251382253; K=3; M=3.864911333333334;  N=11.594734; O=0.2518587252793198; P=0.2518587252793198; Q=5.441912392977112;  R=48.977211536794;   S=0.08701057203081151; T=0.08701057203081153 }
    4  = @{ E=3; G=1.408030333333333;  H=4.224091;  I=0.3454737251382253; J=0.3454737251382253; K=3; M=6.567164333333333;  N=19.701493; O=0.4279522853287917; P=0.4279522853287916; Q=9.246766585318111;  R=83.22089926786299; S=0.1478462701939544;  T=0.1478462701939544 }
    5  = @{ E=3; G=2.015377;           H=6.046131;  I=0.494492045565236;  J=0.4944920455652361; K=3; M=4.913477;           N=14.740431; O=0.3201889893918886; P=0.3201889893918885; Q=9.902508535829;     R=89.122576822461;   S=0.1583309083318606;  T=0.1583309083318606 }
    6  = @{ E=3; G=2.015377;           H=6.046131;  I=0.494492045565236;  J=0.4944920455652361; K=3; M=3.864911333333334;  N=11.594734; O=0.2518587252793198; P=0.2518587252793198; Q=7.789253408239333;  R=70.103280674154;   S=0.1245421362568237;  T=0.1245421362568237 }
    7  = @{ E=3; G=2.015377;           H=6.046131;  I=0.494492045565236;  J=0.4944920455652361; K=3; M=6.567164333333333;  N=19.701493; O=0.4279522853287917; P=0.4279522853287916; Q=13.23531195262033;  R=119.117807573583;  S=0.2116190009765518;  T=0.2116190009765518 }
    8  = @{ E=3; G=0.6522436666666667; H=1.956731;  I=0.1600342292965385; J=0.1600342292965385; K=3; M=4.913477;           N=14.740431; O=0.3201889893918886; P=0.3201889893918885; Q=3.204784254562334;  R=28.843058291061;   S=0.05124119814656844; T=0.05124119814656843 }
    9  = @{ E=3; G=0.6522436666666667; H=1.956731;  I=0.1600342292965385; J=0.1600342292965385; K=3; M=3.864911333333334;  N=11.594734; O=0.2518587252793198; P=0.2518587252793198; Q=2.520863939394889;  R=22.687775454554;   S=0.04030601699168457; T=0.04030601699168457 }
    10 = @{ E=3; G=0.6522436666666667; H=1.956731;  I=0.1600342292965385; J=0.1600342292965385; K=3; M=6.567164333333333;  N=19.701493; O=0.4279522853287917; P=0.4279522853287916; Q=4.283391344375889;  R=38.550522099383;   S=0.06848701415828554; T=0.06848701415828554 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
